$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.388.84"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.508.87"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "311.44"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "98.21"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").Value = "0.560"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "34.93"
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "7.15"
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").Value = "2.896.81"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "15.32"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "2.527.54"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "0.802"
$ws.Range("E17").Value = "  -4.30%  "
$ws.Range("D18").Value = "42.376.26"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "6.55"
$ws.Range("E19").Value = "  -4.68%  "
$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "11.94"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "68.37"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "239.80"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").Value = "2.83"
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "25.19"
$ws.Range("E27").Value = "  -5.52%  "
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "9.91"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").Value = "37.48"
$ws.Range("E30").Value = "  -8.63%  "
$ws.Range("D31").Value = "5.78"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "156.29"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "2.75"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "0.0777"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("D37").Value = "1.94"
$ws.Range("E37").Value = "  -6.68%  "
$ws.Range("D38").Value = "17.09"
$ws.Range("E38").Value = "  -7.51%  "
$ws.Range("D39").Value = "0.106"
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").Value = "21.07"
$ws.Range("E42").Value = "  -6.27%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0293"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.996.20"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "3.21"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").Value = "9.00"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "2.751.78"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "78.19"
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.186"
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("D51").Value = "99.32"
$ws.Range("E51").Value = "  -2.71%  "
